$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the "age" column (D) to the end of the designation/bio/work-ex block:
# old order D,E,F,G = age, designation, bio, work ex
# new order D,E,F,G = designation, bio, work ex, age
$ws.Columns("D").Cut()
$ws.Columns("H").Insert()

# Add the new "skills" column (J) with its data
$ws.Range("J1").Value = "skills"
$ws.Range("J2:J4").Value = "React.js, Angular"
$ws.Columns("J").ColumnWidth = 14.44

# Page orientation was switched to portrait
$ws.PageSetup.Orientation = 1

# Update the saved selection/active cell
$ws.Range("H13").Select()
